$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs", "balls", "fours" stats (columns C:E) for Ravi Bishnoi's
# rows 2 and 4 were swapped. Keep the cells text-formatted (as they
# were originally, t="str") while updating their values.

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "6"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "7"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0"
